# loc_hierarchy_list: add Arabic ("ara") and French ("fra") translations
# alongside the existing English ("eng") rows for every hierarchy_level
# (0..5), and re-flag is_active as literal text "TRUE" (Text-formatted)
# instead of an Excel boolean.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# lang_code, hierarchy_level, hierarchy_level_name
$data = @(
    @("eng", 0, "Country"),
    @("ara", 0, "دولة"),
    @("fra", 0, "Pays"),
    @("eng", 1, "Region"),
    @("ara", 1, "المناطق"),
    @("fra", 1, "Région"),
    @("eng", 2, "Province"),
    @("ara", 2, "المحافظة"),
    @("fra", 2, "Province"),
    @("eng", 3, "City"),
    @("ara", 3, "مدينة"),
    @("fra", 3, "Ville"),
    @("eng", 4, "Zone"),
    @("ara", 4, "منطقة"),
    @("fra", 4, "Zone"),
    @("eng", 5, "Postal Code"),
    @("ara", 5, "رمز بريدي"),
    @("fra", 5, "code postal")
)

$row = 2
foreach ($item in $data) {
    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)
    $dCell = $ws.Cells.Item($row, 4)

    # A/B/C revert to the plain "Normal" style (no explicit xf, matches
    # how the original s="2" formatting was dropped in the update).
    $aCell.Style = "Normal"
    $bCell.Style = "Normal"
    $cCell.Style = "Normal"

    $aCell.Value = $item[0]
    $bCell.Value = $item[1]
    $cCell.Value = $item[2]

    # is_active -> literal text "TRUE" (Text number format) rather than
    # a boolean. Going through a formula + paste-values round trip keeps
    # Excel from re-coercing the literal "TRUE" string back into a bool.
    $dCell.NumberFormat = "@"
    $dCell.Formula = '=TEXT(TRUE,"@")'
    $dCell.Copy()
    $dCell.PasteSpecial(-4163)

    $row = $row + 1
}

$excel.CutCopyMode = $false

$ws.Range("A2:D19").Select()
